$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("SFIA Level") to hold "Skill Description"
$ws.Range("B1").EntireColumn.Insert()

# Header for the new column
$ws.Range("B1").Value = "Skill Description"

# Populate the new column with the same value as column A (SkillCode) for each data row
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 1).Value2
}
